# Add 2022-Q3 data:
#  1. Insert a new worksheet "2022-Q3" right after "总计", populated with
#     the quarter's fund-holdings table (same layout as the other quarter
#     sheets).
#  2. Insert a new row at the top of the "总计" (totals) sheet summarizing
#     the new quarter, pushing the existing rows down by one.
#
# NOTE: worksheet object refs captured via Worksheets.Item(i) are
# position-bound in this host, so any ref grabbed *before* a
# Worksheets.Add() can silently repoint to a different sheet afterwards.
# To stay safe, the template sheet is re-fetched by position *after* the
# new sheet has been inserted, and the clipboard Copy/PasteSpecial pair
# is kept uninterrupted by any sheet-collection change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet right after "总计" (position 2).
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Re-fetch the (still unrenamed) "2022-Q2" sheet, now pushed to position 3,
# and use it as a formatting template so borders/bold header style/index
# -column style match the existing quarter sheets exactly.
$templateSheet = $wb.Worksheets.Item(3)
$templateSheet.Range("A1:H2").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "008372"
$newSheet.Range("B2").Style = "Normal"

$newSheet.Range("C2").Value = "富国阿尔法两年持有期混合"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "8.64"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "72.16"
$newSheet.Range("E2").Style = "Normal"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "3.28"
$newSheet.Range("F2").Style = "Normal"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.2834"
$newSheet.Range("G2").Style = "Normal"

$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 2) Insert a new row 2 in "总计" for the 2022-Q3 summary, shifting the
#    previously existing rows (2022-Q2 ... 2020-Q4) down by one.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Row-2 insert copies row-1's (header) formatting onto the new row; match
# the plain (unstyled) data rows below by copying the index-column style
# from A3 onto A2 and clearing the spurious header style from B2:D2.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.28

Write-Output "2022-Q3 sheet inserted and total sheet updated"
